$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Trend_instructions")

$ws.Range("J1").Value = "override_normalization"
$ws.Range("J2").Value = "T"

$ws.Range("J2").Select()
